{"js": "// Apply the built-in \"Title\" paragraph style to the document's first\n// paragraph (the diff adds <w:pPr><w:pStyle w:val=\"Title\"/></w:pPr> to\n// the sole paragraph in document.xml).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.style = \"Title\";\nawait context.sync();\n", "ps1": "# Apply the built-in \"Title\" paragraph style to the document's first\n# paragraph (the diff adds <w:pPr><w:pStyle w:val=\"Title\"/></w:pPr> to\n# the sole paragraph in document.xml).\n$d = $word.ActiveDocument\n$d.Paragraphs(1).Style = \"Title\"\n"}
